$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "26.894.84"
$ws.Range("E2").Value = "  +0.42%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.551.80"
$ws.Range("E3").Value = "  +0.36%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.33%  "

# Row 5 - BNB
$ws.Range("D5").Value = "206.34"
$ws.Range("E5").Value = "  +0.95%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  +0.19%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.33%  "

# Row 8 - was Solana, becomes Cardano
$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").Value = "0.247"
$ws.Range("E8").Value = "  +0.78%  "

# Row 9 - was Cardano, becomes Solana
$ws.Range("B9").Value = "Solana"
$ws.Range("C9").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D9").Value = "21.51"
$ws.Range("E9").Value = "  +0.81%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +0.33%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +0.24%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "1.773.10"
$ws.Range("E12").Value = "  +0.43%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.557.82"
$ws.Range("E13").Value = "  +0.78%  "

# Row 14 - Polkadot
$ws.Range("E14").Value = "  +0.92%  "

# Row 15 - Polygon
$ws.Range("E15").Value = "  +0.93%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "26.904.25"
$ws.Range("E16").Value = "  +0.47%  "

# Row 17 - Litecoin
$ws.Range("E17").Value = "  +1.15%  "

# Row 18 - BitcoinCash
$ws.Range("D18").Value = "213.79"
$ws.Range("E18").Value = "  +0.09%  "

# Row 19 - ShibaInu
$ws.Range("D19").Value = "0.0₃0685"
$ws.Range("E19").Value = "  +0.75%  "

# Row 20 - Chainlink
$ws.Range("D20").Value = "7.23"
$ws.Range("E20").Value = "  -0.25%  "

# Row 21 - Dai
$ws.Range("E21").Value = "  +0.35%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  -0.68%  "

# Row 23 - Avalanche
$ws.Range("E23").Value = "  +1.44%  "

# Row 24 - Toncoin
$ws.Range("D24").Value = "1.96"
$ws.Range("E24").Value = "  -1.84%  "

# Row 25 - Monero
$ws.Range("D25").Value = "152.96"
$ws.Range("E25").Value = "  +0.21%  "

# Row 26 - Cosmos
$ws.Range("D26").Value = "6.66"
$ws.Range("E26").Value = "  +2.55%  "

# Row 27 - EthereumClassic
$ws.Range("D27").Value = "14.86"
$ws.Range("E27").Value = "  +0.43%  "

# Row 28 - BinanceUSD
$ws.Range("E28").Value = "  +0.34%  "

# Row 29 - Stellar
$ws.Range("E29").Value = "  +1.61%  "

# Row 30 - Hedera
$ws.Range("E30").Value = "  -0.40%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  -0.85%  "

# Row 32 - Filecoin
$ws.Range("E32").Value = "  +1.91%  "

# Row 33 - Maker
$ws.Range("D33").Value = "1.368.76"
$ws.Range("E33").Value = "  +1.02%  "

# Row 34 - InternetComputer(DFINITY)
$ws.Range("E34").Value = "  +2.09%  "

# Row 35 - LidoDAOToken
$ws.Range("E35").Value = "  +3.75%  "

# Row 36 - TrustWalletToken
$ws.Range("E36").Value = "  +6.84%  "

# Row 37 - HuobiToken
$ws.Range("E37").Value = "  +0.47%  "

# Row 38 - VeChain
$ws.Range("E38").Value = "  +1.05%  "

# Row 39 - ImmutableX
$ws.Range("D39").Value = "0.523"
$ws.Range("E39").Value = "  +0.02%  "

# Row 40 - ARBITRUM
$ws.Range("D40").Value = "0.807"
$ws.Range("E40").Value = "  +0.96%  "

# Row 41 - PaxDollar
$ws.Range("E41").Value = "  +0.30%  "

# Row 42 - WEMIXToken
$ws.Range("E42").Value = "  -0.72%  "

# Row 43 - FraxShare
$ws.Range("E43").Value = "  +0.22%  "

# Row 44 - MXToken
$ws.Range("E44").Value = "  +3.35%  "

# Row 45 - Aave
$ws.Range("D45").Value = "63.52"
$ws.Range("E45").Value = "  +1.22%  "

# Row 46 - RenderToken
$ws.Range("D46").Value = "1.73"
$ws.Range("E46").Value = "  -1.60%  "

# Row 47 - RocketPoolETH
$ws.Range("D47").Value = "1.686.42"
$ws.Range("E47").Value = "  +0.33%  "

# Row 48 - Quant
$ws.Range("D48").Value = "86.13"
$ws.Range("E48").Value = "  +0.48%  "

# Row 49 - Cronos
$ws.Range("E49").Value = "  +0.24%  "

# Row 50 - Algorand
$ws.Range("E50").Value = "  +1.17%  "

# Row 51 - USDD
$ws.Range("E51").Value = "  +0.37%  "
